$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with new month columns, copying style from DE1
$ws.Range("DE1").Copy($ws.Range("DF1"))
$ws.Range("DE1").Copy($ws.Range("DG1"))
$ws.Range("DE1").Copy($ws.Range("DH1"))
$ws.Range("DE1").Copy($ws.Range("DI1"))
$ws.Range("DF1").Value = "2024-09"
$ws.Range("DG1").Value = "2024-10"
$ws.Range("DH1").Value = "2024-11"
$ws.Range("DI1").Value = "2024-12"

# Fill in data for 2024-08 (DE, previously blank placeholder), 2024-09 (DF) and 2024-10 (DG)
$ws.Range("DE2").Value = 5598.074
$ws.Range("DF2").Value = 5853.403
$ws.Range("DG2").Value = 6884.949
$ws.Range("DE3").Value = 580.23
$ws.Range("DF3").Value = 553.946
$ws.Range("DG3").Value = 864.198
$ws.Range("DE4").Value = 10.365
$ws.Range("DF4").Value = 9.464
$ws.Range("DG4").Value = 12.552
$ws.Range("DE5").Value = 225.979
$ws.Range("DF5").Value = 197.723
$ws.Range("DG5").Value = 187.677
$ws.Range("DE6").Value = 4.037
$ws.Range("DF6").Value = 3.378
$ws.Range("DG6").Value = 2.726
$ws.Range("DE7").Value = 44.571
$ws.Range("DF7").Value = 33.974
$ws.Range("DG7").Value = 54.479
$ws.Range("DE8").Value = 0.796
$ws.Range("DF8").Value = 0.58
$ws.Range("DG8").Value = 0.791
$ws.Range("DE9").Value = 134.241
$ws.Range("DF9").Value = 109.58
$ws.Range("DG9").Value = 104.517
$ws.Range("DE10").Value = 2.398
$ws.Range("DF10").Value = 1.872
$ws.Range("DG10").Value = 1.518
$ws.Range("DE11").Value = 133.147
$ws.Range("DF11").Value = 172.103
$ws.Range("DG11").Value = 463.556
$ws.Range("DE12").Value = 2.378
$ws.Range("DF12").Value = 2.94
$ws.Range("DG12").Value = 6.733
$ws.Range("DE13").Value = 42.292
$ws.Range("DF13").Value = 40.565
$ws.Range("DG13").Value = 53.968
$ws.Range("DE14").Value = 0.755
$ws.Range("DF14").Value = 0.6929999999999999
$ws.Range("DG14").Value = 0.784
$ws.Range("DE15").Value = 4188.501
$ws.Range("DF15").Value = 4582.548
$ws.Range("DG15").Value = 5331.12
$ws.Range("DE16").Value = 74.81999999999999
$ws.Range("DF16").Value = 78.289
$ws.Range("DG16").Value = 77.432
$ws.Range("DE17").Value = 3560.84
$ws.Range("DF17").Value = 3576.239
$ws.Range("DG17").Value = 4410.12
$ws.Range("DE18").Value = 63.608
$ws.Range("DF18").Value = 61.097
$ws.Range("DG18").Value = 64.054
$ws.Range("DE19").Value = 2198.536
$ws.Range("DF19").Value = 2342.756
$ws.Range("DG19").Value = 2835.644
$ws.Range("DE20").Value = 39.273
$ws.Range("DF20").Value = 40.024
$ws.Range("DG20").Value = 41.186
$ws.Range("DE21").Value = 1362.305
$ws.Range("DF21").Value = 1233.483
$ws.Range("DG21").Value = 1574.476
$ws.Range("DE22").Value = 24.335
$ws.Range("DF22").Value = 21.073
$ws.Range("DG22").Value = 22.868
$ws.Range("DE23").Value = 0
$ws.Range("DF23").Value = 0
$ws.Range("DG23").Value = 0
$ws.Range("DE24").Value = 0
$ws.Range("DF24").Value = 0
$ws.Range("DG24").Value = 0
$ws.Range("DE25").Value = 0
$ws.Range("DF25").Value = 0
$ws.Range("DG25").Value = 0
$ws.Range("DE26").Value = 0
$ws.Range("DF26").Value = 0
$ws.Range("DG26").Value = 0
$ws.Range("DE27").Value = 401.682
$ws.Range("DF27").Value = 808.587
$ws.Range("DG27").Value = 733.323
$ws.Range("DE28").Value = 7.175
$ws.Range("DF28").Value = 13.814
$ws.Range("DG28").Value = 10.651
$ws.Range("DE29").Value = 401.682
$ws.Range("DF29").Value = 808.587
$ws.Range("DG29").Value = 733.323
$ws.Range("DE30").Value = 7.175
$ws.Range("DF30").Value = 13.814
$ws.Range("DG30").Value = 10.651
$ws.Range("DE31").Value = 0
$ws.Range("DF31").Value = 0
$ws.Range("DG31").Value = 0
$ws.Range("DE32").Value = 0
$ws.Range("DF32").Value = 0
$ws.Range("DG32").Value = 0
$ws.Range("DE33").Value = 0
$ws.Range("DF33").Value = 0
$ws.Range("DG33").Value = 0
$ws.Range("DE34").Value = 0
$ws.Range("DF34").Value = 0
$ws.Range("DG34").Value = 0
$ws.Range("DE35").Value = 0
$ws.Range("DF35").Value = 0
$ws.Range("DG35").Value = 0
$ws.Range("DE36").Value = 0
$ws.Range("DF36").Value = 0
$ws.Range("DG36").Value = 0
$ws.Range("DE37").Value = 0
$ws.Range("DF37").Value = 0
$ws.Range("DG37").Value = 0
$ws.Range("DE38").Value = 0
$ws.Range("DF38").Value = 0
$ws.Range("DG38").Value = 0
$ws.Range("DE39").Value = 0
$ws.Range("DF39").Value = 0
$ws.Range("DG39").Value = 0
$ws.Range("DE40").Value = 0
$ws.Range("DF40").Value = 0
$ws.Range("DG40").Value = 0
$ws.Range("DE41").Value = 0
$ws.Range("DF41").Value = 0
$ws.Range("DG41").Value = 0
$ws.Range("DE42").Value = 0
$ws.Range("DF42").Value = 0
$ws.Range("DG42").Value = 0
$ws.Range("DE43").Value = 1055.322
$ws.Range("DF43").Value = 914.6319999999999
$ws.Range("DG43").Value = 877.309
$ws.Range("DE44").Value = 18.852
$ws.Range("DF44").Value = 15.626
$ws.Range("DG44").Value = 12.742
$ws.Range("DE45").Value = 354.251
$ws.Range("DF45").Value = 356.223
$ws.Range("DG45").Value = 676.521
$ws.Range("DE46").Value = 6.328
$ws.Range("DF46").Value = 6.086
$ws.Range("DG46").Value = 9.826000000000001

Write-Output "done"
